$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Matthew Darby (row 3) gets his hours filled in
$ws.Range("B3").Value = "8h 15m"

# Reflect that B3 is now the active/selected cell, as in the saved file
$ws.Range("B3").Select()
